$d = $word.ActiveDocument

# --- 1. Remove the duplicate title paragraph near the end of the document
#        (bold "Play Dia De Los Muertos Slot Game Free | Pros and Cons").
#        Do this (and the next step) before inserting the new paragraph at the
#        top, so paragraph indices below stay valid and text searches remain
#        unambiguous.
$dupTitlePara = $d.Paragraphs.Item(52)
$dupTitlePara.Range.Delete()

# --- 2. Replace the italic paragraph's text (previously the meta description)
#        with the new feature-image prompt text, keeping its italic run
#        formatting intact.
$newImageText = "Create a feature image for Dia de Los Muertos slot game that showcases the colorful and festive atmosphere of the Mexican holiday while featuring a cartoon-style Maya warrior wearing glasses and looking happy. The warrior should be holding a skull or a flower that represents the essence of the game. The background should be a desert setting with a setting sun or a colorful sky to represent the festival of Dia de Los Muertos. The image should make the viewer feel excited and interested in playing the game."
$italicPara = $d.Paragraphs.Item(52)
$italicPara.Range.Find.Execute("Discover the pros and cons of playing Dia De Los Muertos slot game, a colorful and festive game with smooth gameplay and a gamble feature. Play for free!", $true, $false, $false, $false, $false, $true, 1, $false, $newImageText, 2)

# --- 3. Insert a new "Meta description" paragraph right after the document
#        title (Heading1) paragraph.
$p1 = $d.Paragraphs.Item(1)
$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs.Item(2)
$metaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t xml:space="preserve">: Discover the pros and cons of playing Dia De Los Muertos slot game, a colorful and festive game with smooth gameplay and a gamble feature. Play for free!</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p2.Range.InsertXML($metaXml)
